# Update the dMRI tracts table (rows 47-64 -> 47-67):
#  - section header row 47 keeps its text but gets a bottom rule + taller row
#  - a new "fxcut" row is inserted right after "fx"
#  - abbreviation column (A) is re-ordered/re-labelled (lower-case, no tabs)
#  - description column (B) gets new wording, a new "scs"/"sifc" split, etc.
#  - column B is widened to fit the longer descriptions
#  - A48:A67 become non-italic Arial 11, vertically centred + wrapped, with a
#    medium top rule over the first data row under the header

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Make room: one new row is needed (18 rows of data -> 19 after the split
#    of "fx" into "fx" + "fxcut"), so insert a blank row right after row 48.
# ---------------------------------------------------------------------------
$ws.Rows("49:49").Insert()

# ---------------------------------------------------------------------------
# 2. Row 47 - section header ("dMRI tracts"); add the thicker rule under it.
# ---------------------------------------------------------------------------
$ws.Range("A47").Value = "dMRI tracts"
$ws.Range("B47").ClearContents()
$ws.Rows("47:47").RowHeight = 17
$ws.Range("A47:B47").Borders.Item(9).LineStyle = 1
$ws.Range("A47:B47").Borders.Item(9).Weight = 3

# ---------------------------------------------------------------------------
# 3. Data rows 48-67: abbreviation (A) + full description (B).
# ---------------------------------------------------------------------------
$ws.Range("A48").Value = "fx"
$ws.Range("B48").Value = "fornix"
$ws.Range("A49").Value = "fxcut"
$ws.Range("B49").Value = "Fornix excluding fimbria "
$ws.Range("A50").Value = "cgc"
$ws.Range("B50").Value = "cingulate cingulum"
$ws.Range("A51").Value = "cgh"
$ws.Range("B51").Value = "parahippocampal cingulum"
$ws.Range("A52").Value = "cst"
$ws.Range("B52").Value = "corticospinal tract (pyramidal tract)"
$ws.Range("A53").Value = "atr"
$ws.Range("B53").Value = "anterior thalamic radiations"
$ws.Range("A54").Value = "unc"
$ws.Range("B54").Value = "uncinate"
$ws.Range("A55").Value = "ilf"
$ws.Range("B55").Value = "inferior longitudinal fasciculus"
$ws.Range("A56").Value = "ifo"
$ws.Range("B56").Value = "inferior frontal occipital fasciculus"
$ws.Range("A57").Value = "fmaj"
$ws.Range("B57").Value = "forceps major"
$ws.Range("A58").Value = "fmin"
$ws.Range("B58").Value = "forceps minor"
$ws.Range("A59").Value = "cc"
$ws.Range("B59").Value = "corpus callosum"
$ws.Range("A60").Value = "slf"
$ws.Range("B60").Value = "superior longitudinal fasciculus"
$ws.Range("A61").Value = "tslf"
$ws.Range("B61").Value = "temporal superior longitudinal fasciculus (arcuate fasciculus)"
$ws.Range("A62").Value = "pslf"
$ws.Range("B62").Value = "parietal superior longitudinal fasciculus"
$ws.Range("A63").Value = "scs"
$ws.Range("B63").Value = "superior corticostriate"
$ws.Range("A64").Value = "fscs"
$ws.Range("B64").Value = "frontal superior corticostriate"
$ws.Range("A65").Value = "pscs"
$ws.Range("B65").Value = "parietal superior corticostriate"
$ws.Range("A66").Value = "sifc"
$ws.Range("B66").Value = "striatal inferior frontal cortex tract"
$ws.Range("A67").Value = "ifsfc"
$ws.Range("B67").Value = "inferior frontal to superior frontal cortical tract"

# ---------------------------------------------------------------------------
# 4. Formatting for the A48:A67 abbreviation column: Arial 11, not bold/
#    italic, vertically centred + wrapped; a medium rule on top of A48 only.
# ---------------------------------------------------------------------------
$abbrev = $ws.Range("A48:A67")
$abbrev.Font.Name = "Arial"
$abbrev.Font.Size = 11
$abbrev.Font.Bold = $false
$abbrev.Font.Italic = $false
$abbrev.VerticalAlignment = -4108
$abbrev.WrapText = $true

$ws.Range("A48").Borders.Item(8).LineStyle = 1
$ws.Range("A48").Borders.Item(8).Weight = -4138

# ---------------------------------------------------------------------------
# 5. Column B needs to be wider to fit the longer descriptions.
# ---------------------------------------------------------------------------
$ws.Columns("B:B").ColumnWidth = 54.33203125

# ---------------------------------------------------------------------------
# 6. Scroll/selection bookkeeping to match where the author ended up.
# ---------------------------------------------------------------------------
$ws.Application.Goto($ws.Range("A69"))
$ws.Range("A69").Select()
